# Updates cryptos list data (prices / 1h volume %, and a reshuffle of the
# rows for the 43-47 rank block) per the "Updated cryptos list" commit.
#
# Every new value is written with a leading apostrophe, Excel's classic
# "force text" prefix. Without it, number-looking strings such as
# "528.04", "1.00" or "0.0703" get auto-coerced to floating point values
# (losing trailing zeros / switching to scientific notation), which the
# source workbook never does -- every cell in this sheet is stored as
# text. The leading apostrophe itself is a marker Excel strips from the
# stored value, so the cell content ends up exactly as intended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''58.008.59'
$ws.Range('E2').Value = '''  +2.28%  '
$ws.Range('D3').Value = '''3.062.82'
$ws.Range('E3').Value = '''  +2.50%  '
$ws.Range('E4').Value = '''  +0.00%  '
$ws.Range('D5').Value = '''528.04'
$ws.Range('E5').Value = '''  +6.00%  '
$ws.Range('D6').Value = '''143.64'
$ws.Range('E6').Value = '''  +6.34%  '
$ws.Range('E7').Value = '''  -0.02%  '
$ws.Range('E8').Value = '''  +5.16%  '
$ws.Range('D9').Value = '''7.64'
$ws.Range('E9').Value = '''  +5.15%  '
$ws.Range('E10').Value = '''  +6.85%  '
$ws.Range('E11').Value = '''  +5.73%  '
$ws.Range('E12').Value = '''  +2.01%  '
$ws.Range('D13').Value = '''3.588.87'
$ws.Range('E13').Value = '''  +2.68%  '
$ws.Range('D14').Value = '''27.39'
$ws.Range('E14').Value = '''  +7.81%  '
$ws.Range('E15').Value = '''  +16.14%  '
$ws.Range('D16').Value = '''57.963.91'
$ws.Range('E16').Value = '''  +2.29%  '
$ws.Range('D17').Value = '''6.22'
$ws.Range('E17').Value = '''  +8.47%  '
$ws.Range('D18').Value = '''3.064.95'
$ws.Range('E18').Value = '''  +2.76%  '
$ws.Range('D19').Value = '''13.23'
$ws.Range('E19').Value = '''  +7.04%  '
$ws.Range('D20').Value = '''8.23'
$ws.Range('E20').Value = '''  +5.53%  '
$ws.Range('D21').Value = '''342.32'
$ws.Range('E21').Value = '''  +4.78%  '
$ws.Range('E22').Value = '''  -0.02%  '
$ws.Range('D23').Value = '''5.66'
$ws.Range('E23').Value = '''  -1.35%  '
$ws.Range('E24').Value = '''  +7.48%  '
$ws.Range('D25').Value = '''65.44'
$ws.Range('E25').Value = '''  +5.98%  '
$ws.Range('D26').Value = '''0.0₃0983'
$ws.Range('E26').Value = '''  +9.03%  '
$ws.Range('D27').Value = '''0.171'
$ws.Range('E27').Value = '''  +5.21%  '
$ws.Range('E28').Value = '''  +0.66%  '
$ws.Range('E29').Value = '''  +8.94%  '
$ws.Range('E30').Value = '''  +9.82%  '
$ws.Range('E31').Value = '''  +7.09%  '
$ws.Range('E32').Value = '''  +5.65%  '
$ws.Range('D33').Value = '''21.25'
$ws.Range('E33').Value = '''  +3.15%  '
$ws.Range('D34').Value = '''4.83'
$ws.Range('E34').Value = '''  +8.37%  '
$ws.Range('D35').Value = '''157.80'
$ws.Range('E35').Value = '''  +3.56%  '
$ws.Range('E36').Value = '''  +6.95%  '
$ws.Range('D37').Value = '''1.33'
$ws.Range('E37').Value = '''  +4.24%  '
$ws.Range('D38').Value = '''26.24'
$ws.Range('E38').Value = '''  +12.59%  '
$ws.Range('D39').Value = '''0.0703'
$ws.Range('E39').Value = '''  +4.53%  '
$ws.Range('D40').Value = '''3.097.62'
$ws.Range('E40').Value = '''  +2.66%  '
$ws.Range('E41').Value = '''  +3.57%  '
$ws.Range('E42').Value = '''  +11.34%  '
$ws.Range('B43').Value = '''ONDO'
$ws.Range('C43').Value = '''https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D43').Value = '''1.05'
$ws.Range('E43').Value = '''  +4.59%  '
$ws.Range('B44').Value = '''FirstDigitalUSD'
$ws.Range('C44').Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '''1.00'
$ws.Range('E44').Value = '''  +0.18%  '
$ws.Range('B45').Value = '''Mantle'
$ws.Range('C45').Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '''0.667'
$ws.Range('E45').Value = '''  +4.53%  '
$ws.Range('B46').Value = '''Maker'
$ws.Range('C46').Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '''2.347.24'
$ws.Range('E46').Value = '''  +5.68%  '
$ws.Range('B47').Value = '''Stacks'
$ws.Range('C47').Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '''1.48'
$ws.Range('E47').Value = '''  +5.59%  '
$ws.Range('D48').Value = '''6.13'
$ws.Range('E48').Value = '''  +6.71%  '
$ws.Range('E49').Value = '''  +3.75%  '
$ws.Range('E50').Value = '''  +4.34%  '
$ws.Range('D51').Value = '''20.32'
$ws.Range('E51').Value = '''  +6.33%  '
